$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -8.217000000000002
$ws.Range("B9").Value = 5.478999999999999
$ws.Range("D9").Value = -8.187999999999999
$ws.Range("D11").Value = -7.327
$ws.Range("B18").Value = 5.137
$ws.Range("B20").Value = 6.542
$ws.Range("D23").Value = -8.387
$ws.Range("D24").Value = -6.681999999999999
$ws.Range("D26").Value = -7.397
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("D34").Value = -7.685
$ws.Range("B35").Value = 9.35
$ws.Range("D35").Value = -8.013999999999999
$ws.Range("D48").Value = -7.495
$ws.Range("D49").Value = -8.35
$ws.Range("D52").Value = -7.979000000000001
$ws.Range("D66").Value = -7.351000000000001
$ws.Range("D67").Value = -7.419999999999999
$ws.Range("B69").Value = 5.755000000000001
$ws.Range("B76").Value = 6.545
$ws.Range("B78").Value = 8.409000000000001
$ws.Range("D78").Value = -8.141999999999999
$ws.Range("D80").Value = -8.107999999999999
$ws.Range("B82").Value = 5.457
$ws.Range("B83").Value = 5.529999999999999
$ws.Range("B93").Value = 5.575
$ws.Range("D99").Value = -8.253
$ws.Range("D104").Value = -7.647
